# SVM files after running more images
# Adds a new "SVM Output" column (F) with a literal "'N'" text value for
# every data row, widens column E, and updates the active selection —
# matching the authored change while leaving the existing A:E data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "SVM Output"

# New value for every data row (2-11): literal text  'N'
# (typed with a leading quote-prefix so Excel stores it as literal text
# starting with an apostrophe, exactly like the source workbook.)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 6).Value = "''N'"
}

# Column E was widened when the new column was added
$ws.Columns.Item(5).ColumnWidth = 18.75

# Update the selected cell to reflect where the user ended up
$ws.Range("H8").Select()
